$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shared string text updates (header volume number + report week dates)
$ws.Range("A8").Value = "Volume 30   Number  7"
$ws.Range("C9").Value = "Report Covering the Week  2/13/2023  Through  2/19/2023"

# Row 16
$ws.Range("C16").Value = 1
$ws.Range("D16").Value = 2
$ws.Range("E16").Value = -50
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -33.333333333333
$ws.Range("I16").Value = 10
$ws.Range("J16").Value = 15
$ws.Range("K16").Value = -33.333333333333
$ws.Range("L16").Value = 233.333333333333
$ws.Range("M16").Value = -41.176470588235
$ws.Range("N16").Value = -74.358974358974

# Row 17
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 9
$ws.Range("H17").Value = -11.111111111111
$ws.Range("I17").Value = 10
$ws.Range("J17").Value = 13
$ws.Range("K17").Value = -23.076923076923
$ws.Range("L17").Value = 66.666666666666
$ws.Range("M17").Value = -9.090909090909
$ws.Range("N17").Value = -41.176470588235

# Row 18
$ws.Range("C18").Value = 8
$ws.Range("D18").Value = 10
$ws.Range("E18").Value = -20
$ws.Range("F18").Value = 27
$ws.Range("G18").Value = 32
$ws.Range("H18").Value = -15.625
$ws.Range("I18").Value = 53
$ws.Range("J18").Value = 49
$ws.Range("K18").Value = 8.163265306122
$ws.Range("L18").Value = 43.243243243243
$ws.Range("M18").Value = 35.897435897435
$ws.Range("N18").Value = -63.698630136986

# Row 19
$ws.Range("C19").Value = 18
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 157.142857142857
$ws.Range("F19").Value = 71
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = 121.875
$ws.Range("I19").Value = 114
$ws.Range("J19").Value = 70
$ws.Range("K19").Value = 62.857142857142
$ws.Range("L19").Value = 267.741935483871
$ws.Range("M19").Value = 100
$ws.Range("N19").Value = 67.647058823529

# Row 20
$ws.Range("C20").Value = 1
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = -75
$ws.Range("F20").Value = 15
$ws.Range("G20").Value = 9
$ws.Range("H20").Value = 66.666666666666
$ws.Range("I20").Value = 26
$ws.Range("J20").Value = 16
$ws.Range("K20").Value = 62.5
$ws.Range("L20").Value = 116.666666666667
$ws.Range("M20").Value = 52.941176470588
$ws.Range("N20").Value = -93.995381062355

# Row 21
$ws.Range("C21").Value = 31
$ws.Range("D21").Value = 27
$ws.Range("E21").Value = 14.814814814814
$ws.Range("F21").Value = 127
$ws.Range("G21").Value = 91
$ws.Range("H21").Value = 39.560439560439
$ws.Range("I21").Value = 213
$ws.Range("J21").Value = 163
$ws.Range("K21").Value = 30.674846625766
$ws.Range("L21").Value = 139.325842696629
$ws.Range("M21").Value = 51.063829787234
$ws.Range("N21").Value = -69.744318181818

# Row 24
$ws.Range("C24").Value = 9
$ws.Range("D24").Value = 15
$ws.Range("E24").Value = -40
$ws.Range("F24").Value = 47
$ws.Range("G24").Value = 67
$ws.Range("H24").Value = -29.850746268656
$ws.Range("I24").Value = 76
$ws.Range("J24").Value = 117
$ws.Range("K24").Value = -35.042735042735
$ws.Range("L24").Value = 1.333333333333
$ws.Range("M24").Value = -12.643678160919

# Row 25
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 6
$ws.Range("E25").Value = -83.333333333333
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 22
$ws.Range("H25").Value = -63.636363636363
$ws.Range("I25").Value = 25
$ws.Range("J25").Value = 36
$ws.Range("K25").Value = -30.555555555555
$ws.Range("L25").Value = 150
$ws.Range("M25").Value = 19.047619047619

# Row 27
$ws.Range("C27").NumberFormat = "#,##0"
$ws.Range("C27").Value = 1
$ws.Range("D27").NumberFormat = "#,##0"
$ws.Range("D27").Value = 4
$ws.Range("E27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E27").Value = -75
$ws.Range("F27").Value = 1
$ws.Range("G27").NumberFormat = "#,##0"
$ws.Range("G27").Value = 4
$ws.Range("H27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H27").Value = -75
$ws.Range("I27").Value = 2
$ws.Range("J27").NumberFormat = "#,##0"
$ws.Range("J27").Value = 4
$ws.Range("K27").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K27").Value = -50

# Row 30
$ws.Range("D30").NumberFormat = "#,##0"
$ws.Range("D30").Value = 1
$ws.Range("E30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E30").Value = -100
$ws.Range("G30").NumberFormat = "#,##0"
$ws.Range("G30").Value = 1
$ws.Range("H30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("H30").Value = -100
$ws.Range("J30").NumberFormat = "#,##0"
$ws.Range("J30").Value = 1
$ws.Range("K30").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("K30").Value = -100

